$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ip_address_list")
$ws.Activate()

# Helper: write a numeric-looking id as plain text (matches how this sheet
# already stores ids like "514", "51455", ... as text, not numbers)
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 3: becomes the old row 5 data (A=51452422, D=aha added)
Set-TextValue $ws.Range("A3") "51452422"
$ws.Range("D3").Value = "aha"

# Row 4: A becomes old row3's value (5); B becomes old row5's B value (192.168.100.241)
Set-TextValue $ws.Range("A4") "5"
$ws.Range("B4").Value = "192.168.100.241"

# Row 5: becomes old row4 data (A=514, B=192.168.14.241); D5 cleared
Set-TextValue $ws.Range("A5") "514"
$ws.Range("B5").Value = "192.168.14.241"
$ws.Range("D5").ClearContents()

# D1: new, empty text cell (materializes an empty inline-string cell)
$ws.Range("D1").Value = "'"
$ws.Range("D1").Style = "Normal"
